$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.300.32'
$ws.Range("E2").Value = '  -4.93%  '

$ws.Range("D3").Value = '1.563.81'
$ws.Range("E3").Value = '  -4.87%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''1.001'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").Value = '''288.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.73%  '

$ws.Range("D7").Value = '''0.3741'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.19%  '

$ws.Range("D8").Value = '''49.34'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.38%  '

$ws.Range("D9").Value = '''0.3415'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.33%  '

$ws.Range("D10").Value = '''1.164'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.24%  '

$ws.Range("D11").Value = '''0.07638'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.14%  '

$ws.Range("D12").Value = '''1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("E13").Value = '  -3.24%  '

$ws.Range("D14").Value = '''6.013'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.27%  '

$ws.Range("D15").Value = '''6.926'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.35%  '

$ws.Range("D16").Value = '1.566.24'
$ws.Range("E16").Value = '  -4.85%  '

$ws.Range("D17").Value = '''0.00001126'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.78%  '

$ws.Range("D18").Value = '''89.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.66%  '

$ws.Range("D19").Value = '''0.06714'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.78%  '

$ws.Range("D20").Value = '''1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("D21").Value = '''6.226'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.96%  '

$ws.Range("D22").Value = '''16.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.82%  '

$ws.Range("D23").Value = '''0.5287'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.71%  '

$ws.Range("D24").Value = '''11.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.74%  '

$ws.Range("D25").Value = '22.299.14'
$ws.Range("E25").Value = '  -4.92%  '

$ws.Range("D26").Value = '''2.399'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("D27").Value = '''2.793'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.04%  '

$ws.Range("D28").Value = '''20.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.86%  '

$ws.Range("D29").Value = '''145.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.11%  '

$ws.Range("D30").Value = '''4.974'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.06%  '

$ws.Range("D31").Value = '''125.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.70%  '

$ws.Range("D32").Value = '1.735.34'
$ws.Range("E32").Value = '  -5.36%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''1.014'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.50%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''6.175'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.83%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''2.017'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.70%  '

$ws.Range("D36").Value = '''10.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.86%  '

$ws.Range("D37").Value = '''0.08524'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.87%  '

$ws.Range("D38").Value = '''0.02525'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.00%  '

$ws.Range("D39").Value = '''0.2313'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.20%  '

$ws.Range("D40").Value = '''5.497'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.97%  '

$ws.Range("D41").Value = '''1.308'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.16%  '

$ws.Range("D42").Value = '''0.06385'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.91%  '

$ws.Range("D43").Value = '''11.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.06%  '

$ws.Range("D44").Value = '''0.6349'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.74%  '

$ws.Range("D45").Value = '''14.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.39%  '

$ws.Range("D46").Value = '''1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("D47").Value = '''0.5982'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.33%  '

$ws.Range("D48").Value = '''3.745'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.54%  '

$ws.Range("D49").Value = '''2.085'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.91%  '

$ws.Range("D50").Value = '''1.264'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.88%  '

$ws.Range("D51").Value = '''124.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.18%  '
